$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -8
$ws.Range("F12").Value = -7
$ws.Range("F17").Value = -5
$ws.Range("F18").Value = -6
$ws.Range("F19").Value = -1
$ws.Range("F21").Value = -3
$ws.Range("F22").Value = -8
$ws.Range("F24").Value = -5
$ws.Range("F26").Value = -5
$ws.Range("F31").Value = -2
$ws.Range("F32").Value = -6
$ws.Range("F34").Value = 5
$ws.Range("F35").Value = -7
$ws.Range("F36").Value = -3
$ws.Range("F38").Value = -4
$ws.Range("F39").Value = 1
$ws.Range("F43").Value = -2
$ws.Range("F45").Value = -4
$ws.Range("F46").Value = -3
$ws.Range("F49").Value = 8
$ws.Range("F51").Value = 0
$ws.Range("F52").Value = 2
$ws.Range("F53").Value = -1
$ws.Range("F55").Value = 4
